$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Mentioned_in_text")
$ws2 = $wb.Worksheets.Item("Extra_on_github")

# Append a new row (27) to the "Mentioned_in_text" sheet, reusing the
# formatting of row 11 (the "Online Supplementary Material" style block).
$ws1.Range("A11:E11").Copy()
$ws1.Range("A27:E27").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Application.CutCopyMode = $false

$ws1.Range("A27").Value = "Supplementary Figures on GitHub "
$ws1.Range("B27").Value = "Online Supplementary Material"
$ws1.Range("C27").Value = "Metacell Pipeline output figures."
$ws1.Range("D27").Value = "To be prepared"
$ws1.Range("E27").Value = "could be done by re-running clean versions of the metacell scripts."

$ws1.Rows.Item(27).RowHeight = 29

# Update the remembered selections on each sheet.
$null = $ws1.Range("B30").Select()
$null = $ws2.Range("E12").Select()

# Make "Mentioned_in_text" the active (selected) tab, instead of
# "Extra_on_github".
$null = $ws1.Activate()
